$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.133.21"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.669.91"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").Value = "210.65"
$ws.Range("E5").Value = "  -3.91%  "
$ws.Range("D6").Value = "0.5208"
$ws.Range("E6").Value = "  -5.07%  "
$ws.Range("D9").Value = "0.06243"
$ws.Range("E9").Value = "  -3.51%  "
$ws.Range("D10").Value = "21.15"
$ws.Range("E10").Value = "  -3.99%  "
$ws.Range("D11").Value = "0.07516"
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("D12").Value = "1.652.77"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("D13").Value = "4.436"
$ws.Range("E13").Value = "  -2.58%  "
$ws.Range("D14").Value = "0.5588"
$ws.Range("E14").Value = "  -4.46%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.000007978"
$ws.Range("E15").Value = "  -4.89%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "66.22"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "26.174.23"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").Value = "187.00"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("E21").Value = "  -5.76%  "
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").Value = "147.72"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "0.1243"
$ws.Range("E25").Value = "  -6.43%  "
$ws.Range("D26").Value = "7.586"
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("D27").Value = "15.88"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").Value = "0.06191"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").Value = "1.356"
$ws.Range("E29").Value = "  -2.86%  "
$ws.Range("D30").Value = "1.278"
$ws.Range("E30").Value = "  -4.06%  "
$ws.Range("D31").Value = "3.474"
$ws.Range("E31").Value = "  -3.63%  "
$ws.Range("D32").Value = "3.425"
$ws.Range("E32").Value = "  -4.95%  "
$ws.Range("D33").Value = "1.608"
$ws.Range("E33").Value = "  -4.57%  "
$ws.Range("D34").Value = "0.9907"
$ws.Range("E34").Value = "  -5.20%  "
$ws.Range("D35").Value = "0.6026"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D38").Value = "6.126"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").Value = "0.01609"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.072.33"
$ws.Range("E40").Value = "  -4.30%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.8657"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("D43").Value = "99.52"
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("D44").Value = "1.819.06"
$ws.Range("D45").Value = "0.00000000108"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D46").Value = "55.95"
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").Value = "7.944"
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("D50").Value = "0.4252"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").Value = "5.944"
$ws.Range("E51").Value = "  -2.60%  "
